$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Imad"
$ws.Range("C2").Value = "Shehadeh"
$ws.Range("D2").Value = "'25874125"
$ws.Range("E2").Value = "VIP"
